$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) documenting the PF/1.0.4 release:
# column A gets the version label, the rest are marked "X".
$ws.Range("A3").Value = "PF/1.0.4"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
